$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 2
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 2
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 2
